$wb = $excel.ActiveWorkbook

$oldGuid = "db9faf89-8df1-4443-b7e3-12c1f9320384"
$newGuid = "7a771076-b300-48b3-b462-0a68b1bfd62f"

$oldZhHash = "5de3edbfaa1a022becdcc3db6c93c674fbc2bd91"
$newZhHash = "1ff68681f891f401fb14603147f68e6e7bbbdafe"

$oldDeHash = "5de3edbfaa1a022becdcc3db6c93c674fbc2bd91"
$newDeHash = "1ff68681f891f401fb14603147f68e6e7bbbdafe"

# Overview sheet
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("C2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-03-02 10:22:59"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("C2").Value = "$newGuid.$newDeHash.de-de.xlf"
$wsDe.Range("D2").Value = "2016-03-02 10:23:11"
